$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D2").Value = '''69.341.42'
$ws.Range("D2").Style = $plainStyle

$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("D3").Value = '''3.899.85'
$ws.Range("D3").Style = $plainStyle

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("E5").Value = '  +8.57%  '
$ws.Range("D5").Value = '''528.42'
$ws.Range("D5").Style = $plainStyle

$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D6").Value = '''144.91'
$ws.Range("D6").Style = $plainStyle

$ws.Range("E7").Value = '  -1.57%  '
$ws.Range("D7").Value = '''0.613'
$ws.Range("D7").Style = $plainStyle

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("D9").Value = '''0.720'
$ws.Range("D9").Style = $plainStyle

$ws.Range("E10").Value = '  -4.77%  '
$ws.Range("D10").Value = '''0.173'
$ws.Range("D10").Style = $plainStyle

$ws.Range("E11").Value = '  -4.75%  '
$ws.Range("D11").Value = '''0.0000336'
$ws.Range("D11").Style = $plainStyle

$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D12").Value = '''42.03'
$ws.Range("D12").Style = $plainStyle

$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D13").Value = '''4.518.80'
$ws.Range("D13").Style = $plainStyle

$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D14").Value = '''10.27'
$ws.Range("D14").Style = $plainStyle

$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D15").Value = '''3.906.63'
$ws.Range("D15").Style = $plainStyle

$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D16").Value = '''14.01'
$ws.Range("D16").Style = $plainStyle

$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("E17").Value = '  +6.91%  '
$ws.Range("D17").Value = '''1.22'
$ws.Range("D17").Style = $plainStyle

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D18").Value = '''0.134'
$ws.Range("D18").Style = $plainStyle

$ws.Range("E19").Value = '  -1.18%  '
$ws.Range("D19").Value = '''19.76'
$ws.Range("D19").Style = $plainStyle

$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("D20").Value = '''69.277.66'
$ws.Range("D20").Style = $plainStyle

$ws.Range("E21").Value = '  -1.27%  '
$ws.Range("D21").Value = '''425.30'
$ws.Range("D21").Style = $plainStyle

$ws.Range("E22").Value = '  -4.81%  '

$ws.Range("E23").Value = '  -3.89%  '
$ws.Range("D23").Value = '''14.17'
$ws.Range("D23").Style = $plainStyle

$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("D24").Value = '''87.97'
$ws.Range("D24").Style = $plainStyle

$ws.Range("E25").Value = '  +7.72%  '
$ws.Range("D25").Value = '''4.00'
$ws.Range("D25").Style = $plainStyle

$ws.Range("E26").Value = '  -7.47%  '
$ws.Range("D26").Value = '''11.43'
$ws.Range("D26").Style = $plainStyle

$ws.Range("E27").Value = '  -3.56%  '
$ws.Range("D27").Value = '''10.60'
$ws.Range("D27").Style = $plainStyle

$ws.Range("E28").Value = '  -2.26%  '
$ws.Range("D28").Value = '''36.46'
$ws.Range("D28").Style = $plainStyle

$ws.Range("E29").Value = '  -4.31%  '
$ws.Range("D29").Value = '''688.87'
$ws.Range("D29").Style = $plainStyle

$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("D30").Value = '''13.24'
$ws.Range("D30").Style = $plainStyle

$ws.Range("E31").Value = '  -2.75%  '
$ws.Range("D31").Value = '''0.127'
$ws.Range("D31").Style = $plainStyle

$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D32").Value = '''2.85'
$ws.Range("D32").Style = $plainStyle

$ws.Range("E33").Value = '  +10.80%  '
$ws.Range("D33").Value = '''68.26'
$ws.Range("D33").Style = $plainStyle

$ws.Range("B34").Value = 'TheGraph'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("E34").Value = '  +8.06%  '
$ws.Range("D34").Value = '''0.432'
$ws.Range("D34").Style = $plainStyle

$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("D35").Value = '''0.0₃0862'
$ws.Range("D35").Style = $plainStyle

$ws.Range("E36").Value = '  -2.33%  '
$ws.Range("D36").Value = '''5.92'
$ws.Range("D36").Style = $plainStyle

$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("D37").Value = '''39.97'
$ws.Range("D37").Style = $plainStyle

$ws.Range("E38").Value = '  +2.28%  '
$ws.Range("D38").Value = '''0.150'
$ws.Range("D38").Style = $plainStyle

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("E40").Value = '  +8.71%  '
$ws.Range("D40").Value = '''3.34'
$ws.Range("D40").Style = $plainStyle

$ws.Range("E41").Value = '  -0.13%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D42").Value = '''0.0484'
$ws.Range("D42").Style = $plainStyle

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E43").Value = '  +7.17%  '
$ws.Range("D43").Value = '''3.23'
$ws.Range("D43").Style = $plainStyle

$ws.Range("E44").Value = '  -7.31%  '
$ws.Range("D44").Value = '''2.79'
$ws.Range("D44").Style = $plainStyle

$ws.Range("E45").Value = '  +2.08%  '

$ws.Range("E46").Value = '  -1.31%  '

$ws.Range("E47").Value = '  +12.27%  '

$ws.Range("E48").Value = '  +6.46%  '
$ws.Range("D48").Value = '''2.99'
$ws.Range("D48").Style = $plainStyle

$ws.Range("E49").Value = '  +14.74%  '
$ws.Range("D49").Value = '''2.751.55'
$ws.Range("D49").Style = $plainStyle

$ws.Range("E50").Value = '  -8.05%  '
$ws.Range("D50").Value = '''0.0₆0342'
$ws.Range("D50").Style = $plainStyle

$ws.Range("E51").Value = '  +0.08%  '
$ws.Range("D51").Value = '''144.60'
$ws.Range("D51").Style = $plainStyle

